# Applies the cryptos.xlsx data refresh described by the commit:
# "Updated cryptos list on Wed Mar 27 20:13:48 UTC 2024 with GitHub Actions"
#
# Every D/E (and, for three swapped rows, B/C) cell is rewritten with its new
# scraped value. Values are prefixed with a leading apostrophe (the classic
# Excel "store as text" marker) so numeric-looking strings such as "571.70" or
# "0.0000302" stay text cells (matching the original inlineStr/General-format
# cells) instead of being auto-coerced to numbers and losing trailing zeros /
# switching to scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.927.78"
$ws.Range("E2").Value = "'  -1.15%  "
$ws.Range("D3").Value = "'3.503.94"
$ws.Range("E3").Value = "'  -1.86%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'571.70"
$ws.Range("E5").Value = "'  -0.80%  "
$ws.Range("D6").Value = "'183.90"
$ws.Range("E6").Value = "'  -2.56%  "
$ws.Range("D7").Value = "'0.616"
$ws.Range("E7").Value = "'  -2.50%  "
$ws.Range("D8").Value = "'3.498.67"
$ws.Range("E9").Value = "'  +0.13%  "
$ws.Range("D10").Value = "'0.185"
$ws.Range("E10").Value = "'  +4.38%  "
$ws.Range("D11").Value = "'0.647"
$ws.Range("E11").Value = "'  -1.94%  "
$ws.Range("D12").Value = "'54.15"
$ws.Range("E12").Value = "'  -2.59%  "
$ws.Range("D13").Value = "'0.0000302"
$ws.Range("E13").Value = "'  -0.04%  "
$ws.Range("D14").Value = "'9.44"
$ws.Range("E14").Value = "'  -1.77%  "
$ws.Range("D15").Value = "'4.060.60"
$ws.Range("E15").Value = "'  -2.01%  "
$ws.Range("D16").Value = "'19.33"
$ws.Range("E16").Value = "'  -2.10%  "
$ws.Range("D17").Value = "'68.815.54"
$ws.Range("E17").Value = "'  -1.18%  "
$ws.Range("D18").Value = "'3.488.45"
$ws.Range("E18").Value = "'  -2.24%  "
$ws.Range("D19").Value = "'12.25"
$ws.Range("E19").Value = "'  -2.98%  "
$ws.Range("E20").Value = "'  -1.10%  "
$ws.Range("D21").Value = "'544.07"
$ws.Range("E21").Value = "'  +14.71%  "
$ws.Range("E22").Value = "'  -2.53%  "
$ws.Range("D23").Value = "'19.06"
$ws.Range("E23").Value = "'  -1.27%  "
$ws.Range("E24").Value = "'  -0.60%  "
$ws.Range("D25").Value = "'4.39"
$ws.Range("E25").Value = "'  +0.31%  "
$ws.Range("D26").Value = "'94.09"
$ws.Range("E26").Value = "'  -0.62%  "
$ws.Range("D27").Value = "'2.91"
$ws.Range("E27").Value = "'  -3.06%  "
$ws.Range("D28").Value = "'10.78"
$ws.Range("E28").Value = "'  -1.81%  "
$ws.Range("D29").Value = "'9.15"
$ws.Range("E29").Value = "'  -1.96%  "
$ws.Range("D30").Value = "'31.69"
$ws.Range("E30").Value = "'  -2.20%  "
$ws.Range("D31").Value = "'7.23"
$ws.Range("E31").Value = "'  -7.28%  "
$ws.Range("D32").Value = "'12.57"
$ws.Range("E32").Value = "'  +3.16%  "
$ws.Range("D33").Value = "'64.73"
$ws.Range("E33").Value = "'  -2.17%  "
$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "'  -4.61%  "
$ws.Range("D35").Value = "'565.23"
$ws.Range("E35").Value = "'  -3.05%  "
$ws.Range("B36").Value = "'InjectiveProtocol"
$ws.Range("C36").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "'37.84"
$ws.Range("E36").Value = "'  -2.76%  "
$ws.Range("B37").Value = "'Dai"
$ws.Range("C37").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  -0.17%  "
$ws.Range("D38").Value = "'0.396"
$ws.Range("E38").Value = "'  +0.29%  "
$ws.Range("D39").Value = "'2.97"
$ws.Range("E39").Value = "'  +3.97%  "
$ws.Range("D40").Value = "'0.0₃0765"
$ws.Range("E40").Value = "'  -3.83%  "
$ws.Range("D41").Value = "'3.15"
$ws.Range("E41").Value = "'  -1.91%  "
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.132"
$ws.Range("E42").Value = "'  -3.53%  "
$ws.Range("B43").Value = "'Stacks"
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.33"
$ws.Range("E43").Value = "'  -3.16%  "
$ws.Range("D44").Value = "'3.245.89"
$ws.Range("E44").Value = "'  +0.59%  "
$ws.Range("D45").Value = "'3.50"
$ws.Range("E45").Value = "'  +3.65%  "
$ws.Range("D46").Value = "'2.97"
$ws.Range("E46").Value = "'  -3.36%  "
$ws.Range("D47").Value = "'0.0439"
$ws.Range("E47").Value = "'  -0.70%  "
$ws.Range("E48").Value = "'  -2.49%  "
$ws.Range("D49").Value = "'9.01"
$ws.Range("E49").Value = "'  -4.30%  "
$ws.Range("E50").Value = "'  -0.11%  "
$ws.Range("D51").Value = "'139.06"
$ws.Range("E51").Value = "'  +3.61%  "
